$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# Row 70
$ws.Range("H70").Value = 3500
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -10230
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 3500
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -9564
$ws.Range("N73").ClearContents()

# Row 80
$ws.Range("H80").Value = 1262.75
$ws.Range("I80").Value = 419.66666
$ws.Range("K80").Value = 1258.99998
$ws.Range("M80").Value = -260.9999800000001

# Row 83
$ws.Range("H83").Value = 1262.75
$ws.Range("I83").Value = 419.66666
$ws.Range("K83").Value = 3776.99994
$ws.Range("M83").Value = 1215.00006

# Row 86
$ws.Range("H86").Value = 2882.5789
$ws.Range("I86").Value = 3249.4443
$ws.Range("K86").Value = 3249.4443
$ws.Range("M86").Value = -2126.4443

# Row 89
$ws.Range("H89").Value = 2882.5789
$ws.Range("I89").Value = 3249.4443
$ws.Range("K89").Value = 16247.2215
$ws.Range("M89").Value = -10631.2215

# Row 98
$ws.Range("H98").Value = 812.8
$ws.Range("I98").Value = 802.8570999999999
$ws.Range("J98").Value = 836
$ws.Range("K98").Value = 802.8570999999999
$ws.Range("L98").Value = 836
$ws.Range("M98").Value = 695.1429000000001
$ws.Range("N98").Value = -3832

# Row 122
$ws.Range("H122").Value = 812.8
$ws.Range("I122").Value = 802.8570999999999
$ws.Range("J122").Value = 836
$ws.Range("K122").Value = 2408.5713
$ws.Range("L122").Value = 2508
$ws.Range("M122").Value = 41.42870000000039
$ws.Range("N122").Value = -7408

# Row 132
$ws.Range("H132").Value = 4625.591
$ws.Range("I132").Value = 3843.1538
$ws.Range("J132").Value = 5755.778
$ws.Range("K132").Value = 11529.4614
$ws.Range("L132").Value = 17267.334
$ws.Range("M132").Value = -8999.4614
$ws.Range("N132").Value = -22327.334

# Row 138
$ws.Range("H138").Value = 3287.2727
$ws.Range("J138").Value = 5000
$ws.Range("L138").Value = 15000
$ws.Range("N138").Value = -25280

# Row 141
$ws.Range("H141").Value = 1797.7858
$ws.Range("I141").Value = 1859.1538
$ws.Range("K141").Value = 5577.4614
$ws.Range("M141").Value = -397.4614000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -788

# Row 74
$ws.Range("H74").Value = 2239.6
$ws.Range("I74").Value = 2239.6
$ws.Range("K74").Value = 2239.6
$ws.Range("M74").Value = -1365.6

# Row 77
$ws.Range("H77").Value = 2239.6
$ws.Range("I77").Value = 2239.6
$ws.Range("K77").Value = 11198
$ws.Range("M77").Value = -6830

# Row 132
$ws.Range("H132").Value = 866.2727
$ws.Range("I132").Value = 836.55554
$ws.Range("K132").Value = 2509.66662
$ws.Range("M132").Value = 20.33338000000003

# Row 136
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5476.2856
$ws.Range("I86").Value = 4521.8
$ws.Range("K86").Value = 4521.8
$ws.Range("M86").Value = -3398.8

# Row 89
$ws.Range("H89").Value = 5476.2856
$ws.Range("I89").Value = 4521.8
$ws.Range("K89").Value = 22609
$ws.Range("M89").Value = -16993

# Row 134
$ws.Range("H134").Value = 4291.3335
$ws.Range("I134").Value = 4660.4287
$ws.Range("J134").Value = 2999.5
$ws.Range("K134").Value = 13981.2861
$ws.Range("L134").Value = 8998.5
$ws.Range("M134").Value = -11446.2861
$ws.Range("N134").Value = -14068.5

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1588
$ws.Range("I58").Value = 1249.375
$ws.Range("J58").Value = 2942.5
$ws.Range("K58").Value = 1249.375
$ws.Range("L58").Value = 2942.5
$ws.Range("M58").Value = -1046.375
$ws.Range("N58").Value = -3348.5

# Row 132
$ws.Range("H132").Value = 6726.4287
$ws.Range("I132").Value = 6726.4287
$ws.Range("K132").Value = 20179.2861
$ws.Range("M132").Value = -17649.2861

# Row 136
$ws.Range("H136").Value = 1588
$ws.Range("I136").Value = 1249.375
$ws.Range("J136").Value = 2942.5
$ws.Range("K136").Value = 3748.125
$ws.Range("L136").Value = 8827.5
$ws.Range("M136").Value = -1198.125
$ws.Range("N136").Value = -13927.5

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 48.208332
$ws.Range("I2").Value = 29.277779
$ws.Range("J2").Value = 105
$ws.Range("K2").Value = 175.666674
$ws.Range("L2").Value = 630
$ws.Range("M2").Value = -62.666674
$ws.Range("N2").Value = -856

# Row 23
$ws.Range("H23").Value = 336.94736
$ws.Range("I23").Value = 170
$ws.Range("J23").Value = 414
$ws.Range("K23").Value = 510
$ws.Range("L23").Value = 1242
$ws.Range("M23").Value = -275
$ws.Range("N23").Value = -1712

# Row 41
$ws.Range("H41").Value = 450
$ws.Range("I41").Value = 450
$ws.Range("K41").Value = 1350
$ws.Range("M41").Value = -1012

# Row 64
$ws.Range("H64").Value = 1587.1666
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12540

# Row 67
$ws.Range("H67").Value = 1587.1666
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13872

# Row 75
$ws.Range("H75").Value = 1428.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1428.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 4285.5
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -6281.5

# Row 78
$ws.Range("H78").Value = 1428.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1428.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 12856.5
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -22840.5

# Row 81
$ws.Range("I81").Value = 3149.75
$ws.Range("J81").Value = 15
$ws.Range("K81").Value = 9449.25
$ws.Range("L81").Value = 45
$ws.Range("M81").Value = -8326.25
$ws.Range("N81").Value = -2291

# Row 84
$ws.Range("I84").Value = 3149.75
$ws.Range("J84").Value = 15
$ws.Range("K84").Value = 28347.75
$ws.Range("L84").Value = 135
$ws.Range("M84").Value = -22731.75
$ws.Range("N84").Value = -11367

# Row 112
$ws.Range("H112").Value = 35706.23
$ws.Range("I112").Value = 1396.3334
$ws.Range("J112").Value = 45999.2
$ws.Range("K112").Value = 4189.0002
$ws.Range("L112").Value = 137997.6
$ws.Range("M112").Value = -3081.0002
$ws.Range("N112").Value = -140213.6

# Row 129
$ws.Range("H129").Value = 1756.8572
$ws.Range("I129").Value = 1199.1666
$ws.Range("K129").Value = 3597.4998
$ws.Range("M129").Value = 1402.5002

# Row 140
$ws.Range("H140").Value = 1655.4706
$ws.Range("I140").Value = 632.2308
$ws.Range("J140").Value = 4981
$ws.Range("K140").Value = 1896.6924
$ws.Range("L140").Value = 14943
$ws.Range("M140").Value = 3283.3076
$ws.Range("N140").Value = -25303

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 1443500
$ws.Range("I11").Value = 5750000
$ws.Range("K11").Value = 5750000
$ws.Range("M11").Value = -5749861

# Row 80
$ws.Range("H80").Value = 3150.6
$ws.Range("J80").Value = 3418
$ws.Range("L80").Value = 3418
$ws.Range("N80").Value = -5414

# Row 83
$ws.Range("H83").Value = 3150.6
$ws.Range("J83").Value = 3418
$ws.Range("L83").Value = 17090
$ws.Range("N83").Value = -27074

# Row 97
$ws.Range("H97").Value = 966.63635
$ws.Range("I97").Value = 615.3333
$ws.Range("J97").Value = 2547.5
$ws.Range("K97").Value = 615.3333
$ws.Range("L97").Value = 2547.5
$ws.Range("M97").Value = -119.3333
$ws.Range("N97").Value = -3539.5

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 460.66666
$ws.Range("I9").Value = 453
$ws.Range("K9").Value = 453
$ws.Range("M9").Value = -229

# Row 11
$ws.Range("H11").Value = 2603.5
$ws.Range("J11").Value = 2603.5
$ws.Range("L11").Value = 2603.5
$ws.Range("N11").Value = -2883.5

# Row 14
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -2344

# Row 16
$ws.Range("H16").Value = 967.55554
$ws.Range("I16").Value = 838.5
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 838.5
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -668.5
$ws.Range("N16").Value = -2340

# Row 17
$ws.Range("H17").Value = 950
$ws.Range("I17").Value = 950
$ws.Range("K17").Value = 950
$ws.Range("M17").Value = -780

# Row 19
$ws.Range("H19").Value = 2251
$ws.Range("I19").Value = 2003
$ws.Range("K19").Value = 2003
$ws.Range("M19").Value = -1833

# Row 21
$ws.Range("H21").Value = 15000
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15348

# Row 26
$ws.Range("H26").Value = 22000
$ws.Range("I26").Value = 22000
$ws.Range("K26").Value = 22000
$ws.Range("M26").Value = -21705

# Row 30
$ws.Range("H30").Value = 994.4286
$ws.Range("I30").Value = 994.4286
$ws.Range("K30").Value = 994.4286
$ws.Range("M30").Value = -886.4286

# Row 136
$ws.Range("H136").Value = 3627.818
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 423.22223
$ws.Range("J113").Value = 537.25
$ws.Range("L113").Value = 1611.75
$ws.Range("N113").Value = -5951.75

# Row 132
$ws.Range("H132").Value = 2233.3333
$ws.Range("J132").Value = 2233.3333
$ws.Range("L132").Value = 6699.999899999999
$ws.Range("N132").Value = -11759.9999
